$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) values must stay text - force text format, assign, then restore the
# original style so the cell style index in the saved file is unchanged.
$priceCells = @(
    @{ Addr = "D2"; Val = "63.588.07" }
    @{ Addr = "D3"; Val = "3.409.24" }
    @{ Addr = "D5"; Val = "567.77" }
    @{ Addr = "D6"; Val = "157.08" }
    @{ Addr = "D8"; Val = "3.412.79" }
    @{ Addr = "D9"; Val = "0.569" }
    @{ Addr = "D10"; Val = "7.23" }
    @{ Addr = "D11"; Val = "0.118" }
    @{ Addr = "D13"; Val = "3.994.46" }
    @{ Addr = "D15"; Val = "26.93" }
    @{ Addr = "D17"; Val = "63.688.38" }
    @{ Addr = "D18"; Val = "3.404.48" }
    @{ Addr = "D20"; Val = "13.56" }
    @{ Addr = "D21"; Val = "382.84" }
    @{ Addr = "D22"; Val = "7.75" }
    @{ Addr = "D24"; Val = "71.05" }
    @{ Addr = "D25"; Val = "0.514" }
    @{ Addr = "D26"; Val = "0.0000113" }
    @{ Addr = "D27"; Val = "9.69" }
    @{ Addr = "D29"; Val = "0.998" }
    @{ Addr = "D30"; Val = "6.05" }
    @{ Addr = "D31"; Val = "1.38" }
    @{ Addr = "D34"; Val = "22.86" }
    @{ Addr = "D36"; Val = "1.51" }
    @{ Addr = "D37"; Val = "160.41" }
    @{ Addr = "D38"; Val = "0.841" }
    @{ Addr = "D39"; Val = "1.81" }
    @{ Addr = "D40"; Val = "2.820.62" }
    @{ Addr = "D41"; Val = "25.90" }
    @{ Addr = "D42"; Val = "42.97" }
    @{ Addr = "D43"; Val = "0.0718" }
    @{ Addr = "D44"; Val = "6.37" }
    @{ Addr = "D45"; Val = "25.57" }
    @{ Addr = "D46"; Val = "4.34" }
    @{ Addr = "D47"; Val = "0.0302" }
    @{ Addr = "D48"; Val = "328.09" }
)

foreach ($item in $priceCells) {
    $rng = $ws.Range($item.Addr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $item.Val
    $rng.Style = $origStyle
}

# Column E (Volume/1h) values already contain non-numeric characters (%, spaces)
# so plain assignment keeps them as text.
$ws.Range("E2").Value = "  -1.67%  "
$ws.Range("E3").Value = "  -0.34%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("E5").Value = "  -0.84%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -0.43%  "
$ws.Range("E9").Value = "  -7.60%  "
$ws.Range("E10").Value = "  +0.94%  "
$ws.Range("E11").Value = "  -3.62%  "
$ws.Range("E12").Value = "  -4.49%  "
$ws.Range("E13").Value = "  -0.44%  "
$ws.Range("E14").Value = "  +0.10%  "
$ws.Range("E15").Value = "  -3.60%  "
$ws.Range("E16").Value = "  -9.33%  "
$ws.Range("E17").Value = "  -1.49%  "
$ws.Range("E18").Value = "  -0.13%  "
$ws.Range("E19").Value = "  -4.49%  "
$ws.Range("E20").Value = "  -3.10%  "
$ws.Range("E21").Value = "  +1.46%  "
$ws.Range("E22").Value = "  -3.78%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("E24").Value = "  -1.92%  "
$ws.Range("E25").Value = "  -6.79%  "
$ws.Range("E26").Value = "  -5.30%  "
$ws.Range("E27").Value = "  -5.48%  "
$ws.Range("E28").Value = "  +0.14%  "
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("E30").Value = "  -2.61%  "
$ws.Range("E31").Value = "  -7.26%  "
$ws.Range("E32").Value = "  -2.63%  "
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("E34").Value = "  -1.06%  "
$ws.Range("E35").Value = "  -4.15%  "
$ws.Range("E36").Value = "  -6.12%  "
$ws.Range("E37").Value = "  +0.18%  "
$ws.Range("E38").Value = "  +9.05%  "
$ws.Range("E39").Value = "  -4.80%  "
$ws.Range("E40").Value = "  -2.69%  "
$ws.Range("E41").Value = "  -3.15%  "
$ws.Range("E42").Value = "  +0.54%  "
$ws.Range("E43").Value = "  -5.64%  "
$ws.Range("E44").Value = "  -9.61%  "
$ws.Range("E45").Value = "  -3.95%  "
$ws.Range("E46").Value = "  -6.08%  "
$ws.Range("E47").Value = "  -3.94%  "
$ws.Range("E48").Value = "  +1.78%  "
$ws.Range("E49").Value = "  +6.65%  "
$ws.Range("E51").Value = "  -5.63%  "
